# XF-970 AUTO_TC 6.2.7 Read Users listed inside a Tenant - User listed read into a tenant
#
# Adds new parameter columns to the "6_Tenants" sheet: a "Role" block
# (Role Name / Times Applied / Last Edit / no. of permissions) and a
# "Users listed inside a Tenant" block (FIRST NAME / LAST NAME / EMAIL /
# TENANT ADMINISTRATOR / INACTIVE / ACTIONS), each with a matching
# "...Title" parameter-name row above it (row 1 = parameter name,
# row 2 = the actual display text used by the automation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6_Tenants")
$ws.Activate()

# --- "Users listed inside a Tenant" table headers/values (columns AS:AX) ---
$ws.Range("AS2").Value = "FIRST NAME"
$ws.Range("AT2").Value = "LAST NAME"
$ws.Range("AU2").Value = "EMAIL"
$ws.Range("AV2").Value = "TENANT ADMINISTRATOR"
$ws.Range("AW2").Value = "INACTIVE"
$ws.Range("AX2").Value = "ACTIONS"

$ws.Range("AS1").Value = "FirstNameTitle"
$ws.Range("AT1").Value = "LastNameTitle"
$ws.Range("AU1").Value = "EmailTitle"
$ws.Range("AV1").Value = "TenantAdministratorTitle"
$ws.Range("AW1").Value = "InactiveTitle"
$ws.Range("AX1").Value = "ActionsTitle"

# --- "Role" table headers/values (columns AO:AR) ---
$ws.Range("AO2").Value = "Role Name"
$ws.Range("AO1").Value = "Role NameTitle"

$ws.Range("AP1").Value = "TimesAppliedTitle"
$ws.Range("AP2").Value = "Times Applied"

$ws.Range("AQ1").Value = "LastEditTitle"
$ws.Range("AQ2").Value = "Last Edit"

$ws.Range("AR1").Value = "sNoOfPermissionsTitle"
$ws.Range("AR2").Value = "no. of permissions"

# Widen the newly added Role columns to match the rest of the parameter table.
$ws.Range("AO1:AR1").ColumnWidth = 13.1640625

# Extend the blank placeholder cells on rows 2 and 3 the same way the rest
# of the table already does, so the used range grows consistently.
$ws.Range("BM2:BP2").NumberFormat = "@"
$ws.Range("AS3:AY3").NumberFormat = "@"

# Leave the sheet positioned/selected where the author ended up after
# adding the new columns.
$excel.ActiveWindow.ScrollColumn = 34
$ws.Range("AR5").Select()

# A page-setup pass (e.g. Print Preview / Page Setup dialog) was touched
# while the sheet was active, pinning the sheet to portrait orientation.
$ws.PageSetup.Orientation = 1

$wb.Save()
